$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.640.21"
$ws.Range("E2").Value = "  +0.60%  "

$ws.Range("D3").Value = "'1.883.83"
$ws.Range("E3").Value = "  +0.31%  "

$ws.Range("D4").Value = "'1.0000"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'249.57"
$ws.Range("E5").Value = "  +0.99%  "

$ws.Range("D6").Value = "'1.0000"
$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").Value = "'0.4759"
$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").Value = "'0.2941"
$ws.Range("E8").Value = "  +1.39%  "

$ws.Range("D9").Value = "'0.06543"
$ws.Range("E9").Value = "  +0.35%  "

$ws.Range("D10").Value = "'22.05"
$ws.Range("E10").Value = "  +0.92%  "

$ws.Range("D11").Value = "'0.07738"
$ws.Range("E11").Value = "  +0.16%  "

$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").Value = "'0.7405"
$ws.Range("E12").Value = "  +0.36%  "

$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").Value = "'96.86"
$ws.Range("E13").Value = "  -0.34%  "

$ws.Range("D14").Value = "'1.883.93"
$ws.Range("E14").Value = "  +0.32%  "

$ws.Range("E15").Value = "  +2.29%  "

$ws.Range("D16").Value = "'275.13"
$ws.Range("E16").Value = "  +0.79%  "

$ws.Range("D17").Value = "'30.615.22"
$ws.Range("E17").Value = "  +0.57%  "

$ws.Range("D18").Value = "'13.19"
$ws.Range("E18").Value = "  -2.96%  "

$ws.Range("E19").Value = "  -0.41%  "

$ws.Range("D20").Value = "'1.0000"
$ws.Range("E20").Value = "  -0.01%  "

$ws.Range("D21").Value = "'2.129.30"
$ws.Range("E21").Value = "  +0.38%  "

$ws.Range("D22").Value = "'5.351"
$ws.Range("E22").Value = "  +1.95%  "

$ws.Range("D23").Value = "'0.9990"
$ws.Range("E23").Value = "  -0.11%  "

$ws.Range("D24").Value = "'6.239"
$ws.Range("E24").Value = "  +0.97%  "

$ws.Range("D25").Value = "'9.242"
$ws.Range("E25").Value = "  -0.97%  "

$ws.Range("D26").Value = "'164.03"
$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").Value = "'18.87"
$ws.Range("E27").Value = "  +0.16%  "

$ws.Range("D28").Value = "'1.917"
$ws.Range("E28").Value = "  -1.27%  "

$ws.Range("D29").Value = "'1.345"
$ws.Range("E29").Value = "  -1.90%  "

$ws.Range("D30").Value = "'0.09736"
$ws.Range("E30").Value = "  -2.11%  "

$ws.Range("E31").Value = "  -0.99%  "

$ws.Range("D32").Value = "'4.298"
$ws.Range("E32").Value = "  -0.25%  "

$ws.Range("D33").Value = "'4.173"
$ws.Range("E33").Value = "  +2.63%  "

$ws.Range("D34").Value = "'0.04898"
$ws.Range("E34").Value = "  +2.33%  "

$ws.Range("D35").Value = "'1.126"
$ws.Range("E35").Value = "  +0.12%  "

$ws.Range("D36").Value = "'0.7007"
$ws.Range("E36").Value = "  +0.11%  "

$ws.Range("D37").Value = "'2.725"
$ws.Range("E37").Value = "  +0.37%  "

$ws.Range("D38").Value = "'0.01910"
$ws.Range("E38").Value = "  +2.12%  "

$ws.Range("D39").Value = "'2.793"
$ws.Range("E39").Value = "  +2.47%  "

$ws.Range("D40").Value = "'6.311"
$ws.Range("E40").Value = "  -0.51%  "

$ws.Range("D41").Value = "'75.65"
$ws.Range("E41").Value = "  +6.66%  "

$ws.Range("D42").Value = "'2.029"
$ws.Range("E42").Value = "  +4.11%  "

$ws.Range("D43").Value = "'0.4252"
$ws.Range("E43").Value = "  +1.09%  "

$ws.Range("D44").Value = "'0.8427"
$ws.Range("E44").Value = "  +0.74%  "

$ws.Range("D45").Value = "'0.9998"
$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("D46").Value = "'102.73"
$ws.Range("E46").Value = "  -0.10%  "

$ws.Range("E47").Value = "  +1.64%  "

$ws.Range("D48").Value = "'7.065"
$ws.Range("E48").Value = "  -0.26%  "

$ws.Range("D49").Value = "'35.64"
$ws.Range("E49").Value = "  -0.04%  "

$ws.Range("D50").Value = "'917.43"
$ws.Range("E50").Value = "  -0.90%  "

$ws.Range("D51").Value = "'0.05765"
